$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# STEP 1: Insert a new "Meta description: ..." paragraph right after the
# title paragraph ("Play Astro Babes Online Slot for Free - Review").
# ----------------------------------------------------------------------------

# The paragraph near the end of the doc ("Play Astro Babes Online Slot for
# Free - Review", bold) already has the exact run/paragraph-mark shape we
# need (a leading empty run followed by a bold run). Copy it so the new
# paragraph we create matches that shape, then restyle/retext it. (The very
# first paragraph has the same text but is the Heading1 title, so skip it.)
$boldTitlePara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $txt = $cand.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Play Astro Babes Online Slot for Free - Review") {
        $boldTitlePara = $cand
    }
}
if ($boldTitlePara -eq $null) {
    $boldTitlePara = $d.Paragraphs.Item(47)
}
$boldTitlePara.Range.FormattedText.Copy()

$titlePara = $d.Paragraphs.Item(1)
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = $d.Styles.Item("Normal")
$metaRange = $metaPara.Range
$metaRange.FormattedText.Paste()

# Turn the pasted "Play Astro Babes Online Slot for Free - Review" (bold)
# into "Meta description" (still bold), then append the rest of the
# sentence as a normal (non-bold) run.
$metaRange.Find.Execute("Play Astro Babes Online Slot for Free - Review", $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

$metaPara.Range.InsertAfter(": Learn about the unique gameplay style, bonus features and impressive graphics of Astro Babes online slot. Play for free and explore its space theme.")

# ----------------------------------------------------------------------------
# STEP 2: Remove the duplicated bold title paragraph that used to sit right
# before the italic "Learn about..." paragraph near the end of the document,
# and rewrite the italic paragraph's text into the new image-prompt text.
# ----------------------------------------------------------------------------

$dupTitlePara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $txt = $cand.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Play Astro Babes Online Slot for Free - Review") {
        $dupTitlePara = $cand
    }
}
$dupTitlePara.Range.Delete()

$descPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Learn about the unique gameplay style")) {
        $descPara = $cand
    }
}

$descTextRange = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descTextRange.Text = "Prompt: Create a feature image for `"Astro Babes`" that is fitting for the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The image should convey the element of space and the fun and playful vibe of the game. The Maya warrior should be surrounded by Astro Babes and have a space backdrop. Use vibrant colors and playful elements to give a fun and inviting feel to the image."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
